$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.178.63"
$ws.Range("E2").Value = "  +2.78%  "
$ws.Range("D3").Value = "'1.811.45"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'338.83"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").Value = "'0.9996"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "'0.3906"
$ws.Range("E7").Value = "  +2.76%  "
$ws.Range("D8").Value = "'0.3478"
$ws.Range("E8").Value = "  +0.65%  "
$ws.Range("D9").Value = "'48.43"
$ws.Range("E9").Value = "  -0.69%  "
$ws.Range("D10").Value = "'1.189"
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("D11").Value = "'0.07548"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "'22.04"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").Value = "'6.505"
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("D15").Value = "'1.811.13"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").Value = "'7.134"
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").Value = "'0.00001101"
$ws.Range("D18").Value = "'0.06693"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").Value = "'84.92"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "'1.0000"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").Value = "'17.71"
$ws.Range("E21").Value = "  +1.71%  "
$ws.Range("D22").Value = "'6.550"
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "'28.186.66"
$ws.Range("E23").Value = "  +2.76%  "
$ws.Range("D24").Value = "'12.43"
$ws.Range("E24").Value = "  -0.90%  "
$ws.Range("D25").Value = "'2.409"
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("D26").Value = "'1.485"
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("D27").Value = "'2.523"
$ws.Range("E27").Value = "  -1.65%  "
$ws.Range("D28").Value = "'21.27"
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("D29").Value = "'153.83"
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("D30").Value = "'2.020.06"
$ws.Range("E30").Value = "  +1.04%  "
$ws.Range("D31").Value = "'135.58"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("D32").Value = "'6.140"
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("D33").Value = "'4.018"
$ws.Range("E33").Value = "  -0.79%  "
$ws.Range("D34").Value = "'0.08818"
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("D35").Value = "'13.01"
$ws.Range("E35").Value = "  -2.10%  "
$ws.Range("D36").Value = "'0.6940"
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.06548"
$ws.Range("E37").Value = "  +2.46%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "'5.454"
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.02416"
$ws.Range("E39").Value = "  +2.93%  "
$ws.Range("D40").Value = "'1.608"
$ws.Range("E40").Value = "  -2.92%  "
$ws.Range("D41").Value = "'0.2207"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").Value = "'1.254"
$ws.Range("E42").Value = "  -1.55%  "
$ws.Range("D43").Value = "'8.438"
$ws.Range("E43").Value = "  -5.01%  "
$ws.Range("D44").Value = "'14.60"
$ws.Range("E44").Value = "  +0.50%  "
$ws.Range("D45").Value = "'0.6428"
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").Value = "'0.9987"
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "'3.864"
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'2.144"
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'131.52"
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.07197"
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'79.95"
$ws.Range("E51").Value = "  +0.27%  "
